$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" column (C) for rows 2-11: 46064 -> 46065
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 46065
}

# Rows 7-11 had their Beteckning (A), Datum (B) and Area (G) values
# rotated between rows. Apply the new values per row directly.

# Row 7
$ws.Cells.Item(7, 1).Value = "A 57810-2022"
$ws.Cells.Item(7, 2).Value = 44897
$ws.Cells.Item(7, 7).Value = 3.3

# Row 8
$ws.Cells.Item(8, 1).Value = "A 6314-2022"
$ws.Cells.Item(8, 2).Value = 44600
$ws.Cells.Item(8, 7).Value = 3

# Row 9
$ws.Cells.Item(9, 1).Value = "A 33037-2025"
$ws.Cells.Item(9, 2).Value = 45840.39623842593
$ws.Cells.Item(9, 7).Value = 0.8

# Row 10
$ws.Cells.Item(10, 1).Value = "A 33033-2025"
$ws.Cells.Item(10, 2).Value = 45840.39188657407
$ws.Cells.Item(10, 7).Value = 0.7

# Row 11
$ws.Cells.Item(11, 1).Value = "A 25610-2024"
$ws.Cells.Item(11, 2).Value = 45463
$ws.Cells.Item(11, 7).Value = 2.9
